$wb = $excel.ActiveWorkbook

# Rename sheet "06_macros" -> "06_macro"
$ws = $wb.Worksheets.Item("06_macros")
$ws.Name = "06_macro"

# Make sure this sheet is active, then move the selection from C2 to G8
$ws.Activate()
$ws.Range("G8").Select()
